# Apply the trading-journal update:
#  - Row 2 (Greaves Cotton / Long): fill in exit price (G2), mark as a
#    trailing-stop loss, compute L2/S2 P&L, close the trade.
#  - Row 4 (City Union Bank / Long): same pattern.
#  - Row 18 (LLOYDSME / Long): fill in exit price, mark SL Hit, close.
#  - Row 33: brand-new trade row (Patanjali Foods, Long) with its full
#    set of values/formulas.
#  - Move the active selection to K33.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 2: Greaves Cotton ----------------------------------------------
$ws.Range("G2").Value = 204.75
$ws.Range("K2").Value = "Loss"
$ws.Range("L2").Formula = "=(G2-D2)*H2"
$ws.Range("N2").Value = "Trailing Stop Loss Hit"
$ws.Range("T2").Value = "Closed"

# ---- Row 4: City Union Bank ----------------------------------------------
$ws.Range("G4").Value = 2582.8200000000002
$ws.Range("K4").Value = "Loss"
$ws.Range("L4").Formula = "=(G4-D4)*H4"
$ws.Range("N4").Value = "Trailing Stop Loss Hit"
$ws.Range("T4").Value = "Closed"

# ---- Row 18: LLOYDSME -----------------------------------------------------
$ws.Range("G18").Value = 1410.94
$ws.Range("K18").Value = "Loss"
$ws.Range("L18").Formula = "=(G18-D18)*H18"
$ws.Range("N18").Value = "SL Hit"
$ws.Range("T18").Value = "Closed"

# ---- Row 33: new trade (Patanjali Foods) ---------------------------------
# Pull the date format from the cell above (A32) instead of assigning a
# NumberFormat string directly, so the cell reuses the existing "date"
# cell style (s="2") rather than Excel minting a brand-new style.
$ws.Range("A32").Copy()
$ws.Range("A33").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("A33").Value2 = 45870
$ws.Range("B33").Value = "Patanjali Foods"
$ws.Range("C33").Value = "Long"
$ws.Range("D33").Value = 1864.5
$ws.Range("E33").Value = 1814.5
$ws.Range("F33").Value = 2064
$ws.Range("G33").Value = 1813.41
$ws.Range("H33").Value = 60
$ws.Range("I33").Formula = "=(D33-E33)*H33"
$ws.Range("J33").Formula = "=(F33-D33)/(D33-E33)"
$ws.Range("K33").Value = "Loss"
$ws.Range("L33").Formula = "=(G33-D33)*H33"
$ws.Range("M33").Value = "PULL BACK "
$ws.Range("N33").Value = "SL Hit"
$ws.Range("Q33").Value = "Weekly Day"
$ws.Range("R33").Value = "INR"
$ws.Range("S33").Formula = "=L33"
$ws.Range("T33").Value = "Closed"
$ws.Range("U33").Formula = "=H33*D33"

# ---- Selection -------------------------------------------------------------
$ws.Range("K33").Select()

$wb.Save()
